$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Results")

# Mark "Create Test Passed" (column B) as TRUE for all test result rows (2-24),
# completing the create test process results.
$ws.Range("B2:B24").Value = $true
